$d = $word.ActiveDocument

# --- Paragraph 1: "Push" title - shrink sz/szCs from 32 -> 24 (16pt -> 12pt) ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Size = 12
$p1.Range.Font.SizeBi = 12

# --- Paragraph 2: "The git push command ... / It's / the counterpart ..." ---
# add sz/szCs = 24 (12pt) to paragraph mark and all runs
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Font.Size = 12
$p2.Range.Font.SizeBi = 12

# --- Paragraph 3: "You can do this with " -> "You can do this with the following command:" ---
$p3 = $d.Paragraphs.Item(3)
$find = $p3.Range.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "You can do this with "
$find.Replacement.Text = "You can do this with the following command:"
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Font.Size = 12
$p3.Range.Font.SizeBi = 12

# --- Paragraph 4: "$ git push " - switch rFonts from theme (cstheme=minorHAnsi) to
#     explicit Times New Roman for ascii/hAnsi/cs, and add sz/szCs = 24 ---
$p4 = $d.Paragraphs.Item(4)
$xml4 = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">$ git push </w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p4.Range.InsertXML($xml4)

# --- Paragraph 5: inline image + trailing space run - add sz/szCs = 24 ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Font.Size = 12
$p5.Range.Font.SizeBi = 12

# --- Paragraph 6: Wingdings arrow symbol - add sz/szCs = 24 ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Font.Size = 12
$p6.Range.Font.SizeBi = 12

# --- Paragraph 7 (final paragraph, screenshot image): previously had NO pPr at all.
#     Give it an explicit pPr (rFonts Times New Roman + sz/szCs=24) and bump the
#     drawing run's rPr with sz/szCs=24 too. ---
$p7 = $d.Paragraphs.Item(7)
$beforeCount = $d.Paragraphs.Count
$xml7 = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"
             xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"
             xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"
             xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"
             xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"
             xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing">
<w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>
<w:drawing>
<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="598EFD70" wp14:editId="5D9D9D35">
<wp:extent cx="3131820" cy="3005424"/>
<wp:effectExtent l="19050" t="19050" r="11430" b="24130"/>
<wp:docPr id="21" name="Picture 21"/>
<wp:cNvGraphicFramePr>
<a:graphicFrameLocks noChangeAspect="1"/>
</wp:cNvGraphicFramePr>
<a:graphic>
<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">
<pic:pic>
<pic:nvPicPr>
<pic:cNvPr id="1" name=""/>
<pic:cNvPicPr/>
</pic:nvPicPr>
<pic:blipFill>
<a:blip r:embed="rId5"/>
<a:stretch>
<a:fillRect/>
</a:stretch>
</pic:blipFill>
<pic:spPr>
<a:xfrm>
<a:off x="0" y="0"/>
<a:ext cx="3178334" cy="3050061"/>
</a:xfrm>
<a:prstGeom prst="rect">
<a:avLst/>
</a:prstGeom>
<a:ln>
<a:solidFill>
<a:schemeClr val="tx1"/>
</a:solidFill>
</a:ln>
</pic:spPr>
</pic:pic>
</a:graphicData>
</a:graphic>
</wp:inline>
</w:drawing>
</w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p7.Range.InsertXML($xml7)

# InsertXML on the last paragraph of the body leaves behind one extra empty
# trailing paragraph (the original end-of-story mark) - remove it so the
# paragraph count matches the original document again.
if ($d.Paragraphs.Count -gt $beforeCount) {
    $newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
    $cleanupRange = $d.Range($newLast.Range.Start - 1, $newLast.Range.End)
    $cleanupRange.Delete()
}
